$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# "Tra cứu hội viên" (row 13) is now 100% done.
$ws.Range("F13").Value = 1
$ws.Range("F13").NumberFormat = "0%"

# "Xóa hội viên" (row 14) is 90% done, with a note about the issue found.
$ws.Range("F14").Value = 0.9
$ws.Range("F14").NumberFormat = "0%"
$ws.Range("G14").Value = "Một số chuỗi HoTen ko tra cứu được ??. Một số chức năng sẽ bổ sung sau."

# Re-declare all merged ranges so the mergeCells list gets re-emitted in the
# same order Excel produced (A17:A18 first, then the rest in sheet order).
$ws.Range("A17:A18").UnMerge()
$ws.Range("B2:E3").UnMerge()
$ws.Range("A6:A7").UnMerge()
$ws.Range("A8:A10").UnMerge()
$ws.Range("A11:A14").UnMerge()
$ws.Range("A15:A16").UnMerge()

$ws.Range("A17:A18").Merge()
$ws.Range("B2:E3").Merge()
$ws.Range("A6:A7").Merge()
$ws.Range("A8:A10").Merge()
$ws.Range("A11:A14").Merge()
$ws.Range("A15:A16").Merge()

# Move the active selection to F14, matching the author's last edit spot.
$ws.Range("F14").Select() | Out-Null
